$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "59.172.76"
Set-TextValue "E2" "  +6.06%  "

Set-TextValue "D3" "2.543.23"
Set-TextValue "E3" "  +7.09%  "

Set-TextValue "E4" "  +0.43%  "

Set-TextValue "D5" "506.51"
Set-TextValue "E5" "  +6.25%  "

Set-TextValue "D6" "157.60"
Set-TextValue "E6" "  +7.28%  "

Set-TextValue "E7" "  +23.16%  "

Set-TextValue "D8" "0.993"
Set-TextValue "E8" "  -0.62%  "

Set-TextValue "D9" "2.587.82"
Set-TextValue "E9" "  +8.92%  "

Set-TextValue "D10" "6.17"
Set-TextValue "E10" "  +13.74%  "

Set-TextValue "D11" "0.103"
Set-TextValue "E11" "  +6.51%  "

Set-TextValue "D12" "0.341"
Set-TextValue "E12" "  +5.81%  "

Set-TextValue "E13" "  +1.50%  "

Set-TextValue "D14" "2.987.20"
Set-TextValue "E14" "  +7.09%  "

Set-TextValue "D15" "59.074.54"
Set-TextValue "E15" "  +5.72%  "

Set-TextValue "D16" "21.98"
Set-TextValue "E16" "  +8.21%  "

Set-TextValue "E17" "  +4.34%  "

Set-TextValue "D18" "2.581.42"
Set-TextValue "E18" "  +8.54%  "

Set-TextValue "E19" "  +3.13%  "

Set-TextValue "D20" "335.76"
Set-TextValue "E20" "  +6.71%  "

Set-TextValue "D21" "10.38"
Set-TextValue "E21" "  +7.13%  "

Set-TextValue "D22" "6.05"
Set-TextValue "E22" "  +6.78%  "

Set-TextValue "E23" "  +0.68%  "

Set-TextValue "D24" "60.27"
Set-TextValue "E24" "  +6.31%  "

Set-TextValue "E25" "  +5.44%  "

Set-TextValue "E26" "  +7.60%  "

Set-TextValue "D27" "2.669.11"
Set-TextValue "E27" "  +7.27%  "

Set-TextValue "E28" "  -0.63%  "

Set-TextValue "D29" "7.49"
Set-TextValue "E29" "  +3.51%  "

Set-TextValue "D30" "0.0₃0829"
Set-TextValue "E30" "  +8.10%  "

Set-TextValue "D31" "1.00"
Set-TextValue "E31" "  +0.03%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "157.07"
Set-TextValue "E32" "  +7.27%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "19.47"
Set-TextValue "E33" "  +8.15%  "

Set-TextValue "E34" "  +5.86%  "

Set-TextValue "E35" "  +8.73%  "

Set-TextValue "E36" "  +9.73%  "

Set-TextValue "E37" "  +8.05%  "

Set-TextValue "D38" "0.859"
Set-TextValue "E38" "  +2.97%  "

Set-TextValue "D39" "3.79"
Set-TextValue "E39" "  +12.22%  "

Set-TextValue "E40" "  +7.32%  "

Set-TextValue "D41" "291.76"
Set-TextValue "E41" "  +14.99%  "

Set-TextValue "D42" "34.93"
Set-TextValue "E42" "  +4.32%  "

Set-TextValue "E43" "  +7.43%  "

Set-TextValue "E44" "  +7.91%  "

Set-TextValue "D45" "0.0559"
Set-TextValue "E45" "  +4.35%  "

Set-TextValue "D46" "0.990"
Set-TextValue "E46" "  -0.80%  "

Set-TextValue "D47" "19.52"
Set-TextValue "E47" "  +15.46%  "

Set-TextValue "E48" "  +6.87%  "

Set-TextValue "E49" "  +7.07%  "

$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D50" "10.30"
Set-TextValue "E50" "  +1.10%  "

$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D51" "0.718"
Set-TextValue "E51" "  +13.51%  "
